# Apply the gh-pages data refresh (commit 456a3b4) to 合肥-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" ---
$ws1.Range("D4").Value = "滨河西路百大东兴家园东南侧约60米 合肥瑶海富茂大饭店"
$ws1.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202407/vOXvDIpB1720514081756.jpeg"

$ws1.Range("F5").Value = 0
$ws1.Range("F8").Value = 66
$ws1.Range("F12").Value = 107
$ws1.Range("F13").Value = 0
$ws1.Range("F15").Value = 0
$ws1.Range("F19").Value = 5071
$ws1.Range("F22").Value = 0

# --- Sheet "全部类型" ---
$ws4.Range("F2").Value = 0
$ws4.Range("F3").Value = 0

$ws4.Range("D4").Value = "滨河西路百大东兴家园东南侧约60米 合肥瑶海富茂大饭店"
$ws4.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202407/vOXvDIpB1720514081756.jpeg"

$ws4.Range("F8").Value = 66
$ws4.Range("F10").Value = 1290
$ws4.Range("F11").Value = 18
$ws4.Range("F13").Value = 403
$ws4.Range("F16").Value = 401
$ws4.Range("F20").Value = 5071
$ws4.Range("F21").Value = 0
$ws4.Range("F22").Value = 0
$ws4.Range("F25").Value = 212
$ws4.Range("F26").Value = 0

$wb.Save()
